$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Fix capitalization of connector words (de/del/la/el -> De/Del/La/El) in state/municipality names ---
$ws.Range('B12').Value = 'Mazapa De Madero'
$ws.Range('B28').Value = 'Hidalgo Del Parral'
$ws.Range('B38').Value = 'San Francisco De Borja'
$ws.Range('B40').Value = 'Valle De Zaragoza'
$ws.Range('A49').Value = 'Ciudad De México'
$ws.Range('B60').Value = 'Coneto De Comonfort'
$ws.Range('B65').Value = 'Nombre De Dios'
$ws.Range('A72').Value = 'Estado De México'
$ws.Range('B73').Value = 'Ecatepec De Morelos'
$ws.Range('B77').Value = 'Tlalnepantla De Baz'
$ws.Range('B86').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B87').Value = 'Silao De La Victoria'
$ws.Range('B90').Value = 'Acapulco De Juárez'
$ws.Range('B91').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B92').Value = 'Coyuca De Catalán'
$ws.Range('B93').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B99').Value = 'Cuautepec De Hinojosa'
$ws.Range('B102').Value = 'Mixquiahuala De Juárez'
$ws.Range('B104').Value = 'Tulancingo De Bravo'
$ws.Range('B107').Value = 'Autlán De Navarro'
$ws.Range('B112').Value = 'Encarnación De Díaz'
$ws.Range('B116').Value = 'Lagos De Moreno'
$ws.Range('B120').Value = 'San Juan De Los Lagos'
$ws.Range('B122').Value = 'San Miguel El Alto'
$ws.Range('B123').Value = 'San Sebastián Del Oeste'
$ws.Range('B125').Value = 'Tamazula De Gordiano'
$ws.Range('B130').Value = 'Unión De Tula'
$ws.Range('B155').Value = 'Tetela Del Volcán'
$ws.Range('B159').Value = 'Amatlán De Cañas'
$ws.Range('B162').Value = 'Santa María Del Oro'
$ws.Range('B167').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B169').Value = 'Ocotlán De Morelos'
$ws.Range('B170').Value = 'Putla Villa De Guerrero'
$ws.Range('B181').Value = 'Tlacolula De Matamoros'
$ws.Range('B185').Value = 'Cuayuca De Andrade'
$ws.Range('B193').Value = 'Tepexi De Rodríguez'
$ws.Range('B200').Value = 'Amealco De Bonfil'
$ws.Range('B202').Value = 'Pinal De Amoles'
$ws.Range('B203').Value = 'San Juan Del Río'
$ws.Range('B210').Value = 'Santa María Del Río'
$ws.Range('B212').Value = 'Villa De Arista'
$ws.Range('B213').Value = 'Villa De Ramos'
$ws.Range('B241').Value = 'Muñoz De Domingo Arenas'
$ws.Range('B275').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B277').Value = 'Villa De Cos'

# --- Remove trailing metadata/footer rows (284-288) ---
$ws.Range("A284:D288").EntireRow.Delete()

Write-Host "Edit applied successfully"
